$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.380.99"
$ws.Range("E2").Value = "  +3.75%  "
$ws.Range("D3").Value = "'3.246.74"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'577.97"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").Value = "'181.87"
$ws.Range("E6").Value = "  +6.85%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -3.25%  "
$ws.Range("D9").Value = "'3.245.45"
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("E10").Value = "  +5.70%  "
$ws.Range("E11").Value = "  +2.94%  "
$ws.Range("E12").Value = "  +4.98%  "
$ws.Range("D13").Value = "'3.810.91"
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").Value = "'28.43"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").Value = "'67.344.84"
$ws.Range("E16").Value = "  +3.79%  "
$ws.Range("E17").Value = "  +3.03%  "
$ws.Range("D18").Value = "'3.253.00"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").Value = "'5.84"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "'13.50"
$ws.Range("E20").Value = "  +5.18%  "
$ws.Range("D21").Value = "'376.20"
$ws.Range("E21").Value = "  +5.47%  "
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'71.21"
$ws.Range("E24").Value = "  +2.71%  "
$ws.Range("D25").Value = "'0.510"
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").Value = "'9.57"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "'5.77"
$ws.Range("E30").Value = "  +7.99%  "
$ws.Range("E31").Value = "  +3.31%  "
$ws.Range("D32").Value = "'22.64"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'1.27"
$ws.Range("E34").Value = "  +5.33%  "
$ws.Range("D35").Value = "'6.91"
$ws.Range("E35").Value = "  +4.05%  "
$ws.Range("D36").Value = "'163.90"
$ws.Range("E36").Value = "  +5.75%  "
$ws.Range("E37").Value = "  +3.59%  "
$ws.Range("D38").Value = "'0.847"
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("E39").Value = "  +4.42%  "
$ws.Range("D40").Value = "'6.83"
$ws.Range("E40").Value = "  +13.10%  "
$ws.Range("D41").Value = "'26.64"
$ws.Range("E41").Value = "  +2.05%  "
$ws.Range("E42").Value = "  +9.95%  "
$ws.Range("E43").Value = "  +4.81%  "
$ws.Range("D44").Value = "'357.71"
$ws.Range("E44").Value = "  +9.87%  "
$ws.Range("D45").Value = "'2.717.95"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("D46").Value = "'25.58"
$ws.Range("E46").Value = "  +6.03%  "
$ws.Range("D47").Value = "'40.79"
$ws.Range("E47").Value = "  +3.35%  "
$ws.Range("D48").Value = "'0.0678"
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  +6.12%  "
$ws.Range("E51").Value = "  -0.63%  "
